$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.624.16"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.82%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.821.84"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.65%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "305.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.48%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4685"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.75%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3597"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07127"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9021"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07814"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.42"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.791.45"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.257"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.336"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.64%  "
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008561"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.672.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.39%  "
$ws.Range("E22").Value = "  +0.66%  "
$ws.Range("E23").Value = "  +0.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.936"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.973"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "113.62"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.54%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.807"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.96%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08806"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.144"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.774"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7307"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.444"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.123"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.077"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("E37").Value = "  -0.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.920"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05117"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5059"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.94%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.822"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1497"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.997"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4676"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.09%  "
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "99.04"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.95%  "
$ws.Range("E48").Value = "  -1.58%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06009"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "63.78"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "35.74"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.17%  "
